$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 values are written as text (not auto-converted to numbers),
# matching the inlineStr type used in the source file.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.313.74'
$ws.Range('D3').Value = '1.832.23'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '243.24'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '0.6192'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '0.07374'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').Value = '0.2928'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').Value = '23.27'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = '0.07655'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '1.849.25'
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').Value = '4.996'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').Value = '0.6765'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '82.84'
$ws.Range('D16').Value = '0.000008974'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = '5.892'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '29.296.15'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '2.099.39'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').Value = '240.15'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').Value = '7.399'
$ws.Range('E23').Value = '  +2.57%  '
$ws.Range('D24').Value = '0.9995'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').Value = '158.46'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('E26').Value = '  -1.26%  '
$ws.Range('D27').Value = '8.568'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').Value = '17.69'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').Value = '1.492'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').Value = '0.05860'
$ws.Range('E30').Value = '  +4.37%  '
$ws.Range('D31').Value = '1.231'
$ws.Range('E31').Value = '  +2.05%  '
$ws.Range('D32').Value = '4.089'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').Value = '4.103'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').Value = '1.861'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('D36').Value = '0.7217'
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('D37').Value = '2.617'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').Value = '2.860'
$ws.Range('E38').Value = '  +2.89%  '
$ws.Range('D39').Value = '1.222.12'
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('D41').Value = '0.9128'
$ws.Range('E41').Value = '  +1.78%  '
$ws.Range('D42').Value = '6.224'
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = '2.017.34'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('D45').Value = '101.93'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').Value = '65.87'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').Value = '0.5052'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('E48').Value = '  +6.63%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.229'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '0.4056'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').Value = '  -4.42%  '

# Reset style index back to default (no explicit style) now that values are set,
# so cells don't carry a stray text-format style attribute.
$ws.Range("D2:E51").Style = "Normal"
